$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 : БИВТ-22-17
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("БИВТ-22-17")

# F6 was a text "pass" grade placeholder -> becomes the actual numeric score
$ws1.Range("F6").Value = 5

# move the cursor / selection like the author left it
$ws1.Range("N10").Select()

# ---------------------------------------------------------------------------
# Sheet 2 : БИВТ-22-18
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("БИВТ-22-18")

# F6 was blank -> filled in with a real grade
$ws2.Range("F6").Value = 5

$ws2.Range("G20").Select()

# ---------------------------------------------------------------------------
# Sheet 3 : БИВТ-22-20
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("БИВТ-22-20")

# D9 was a text placeholder "5 not verified" -> becomes the real numeric grade
$ws3.Range("D9").Value = 5

# E10 was blank -> filled in with a real grade
$ws3.Range("E10").Value = 4

# F18 was a text placeholder "erased" -> becomes the real numeric grade
$ws3.Range("F18").Value = 5

# Copy formats from a genuinely blank area over a handful of rows (the
# author appears to have used Format Painter / paste-special to tidy up
# stray formatting to the right of the table)
$ws3.Range("R100:Y100").Copy()
$ws3.Range("R3:Y4").PasteSpecial(-4122)
$ws3.Range("Q10:Y10").PasteSpecial(-4122)
$ws3.Range("Q12:Y12").PasteSpecial(-4122)
$ws3.Range("Q15:Y15").PasteSpecial(-4122)
$ws3.Range("Q18:Y18").PasteSpecial(-4122)
$ws3.Range("Q27:Y27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the now-unused scratch rows under the table
$ws3.Rows("32:39").Delete()

# Trim the stray, empty, right-hand column that was part of the used range
$ws3.Columns("AB:AB").Delete()

$ws3.Range("X1").Select()

$ws1.Select()
